# CORE_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer footnote
#  - refresh the Weight / Percent Change figures for rows 2-8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet carries cell-protection; unprotect it so the cells below can be written.
$ws.Unprotect()

# Update the "Model holdings provided as of ..." disclaimer text (A11)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."

# Refresh Weight (D) and Percent Change (E) values
$ws.Range("D2").Value = 0.5024045227354528
$ws.Range("E2").Value = -0.01354491733616625

$ws.Range("D3").Value = 0.2420864354516471
$ws.Range("E3").Value = -0.004268472181336413

$ws.Range("D4").Value = 0.09557288443601987
$ws.Range("E4").Value = -0.008995502248875686

$ws.Range("D5").Value = 0.1036276991769866
$ws.Range("E5").Value = -0.009785093735711015

$ws.Range("D6").Value = 0.02971797913572129
$ws.Range("E6").Value = -0.007079307375872923

$ws.Range("D7").Value = 0.02659047906417222
$ws.Range("E7").Value = -0.004625999686372895

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = -0.01004549004861366
